# Adicionado seleção de substatus
# Duplicate the ICCID/RASTREIO pair from rows 4:6 into new rows 17:19
# (same values/styles as the source rows), then select the newly added
# RASTREIO cells, mirroring the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4:B6").Copy()
$ws.Range("A17").PasteSpecial()

[void]$ws.Range("B17:B19").Select()
